# Add a new "2022-Q4" worksheet (quarterly holdings snapshot) right after
# the "总计" (total) sheet, and record it in the summary table.

$wb = $excel.ActiveWorkbook

$totalWs = $wb.Worksheets.Item("总计")
$q3Ws    = $wb.Worksheets.Item("2022-Q3")

# --- 1. Create the new "2022-Q4" sheet by copying "2022-Q3"'s structure/
#        formatting, placed right after "总计" so every older quarter shifts
#        one tab to the right. -----------------------------------------
$q3Ws.Copy($null, $totalWs)
$q4Ws = $wb.Worksheets.Item(2)
$q4Ws.Name = "2022-Q4"

# Update the fund rows with the new quarter's figures (fund codes / names
# stay the same as they were copied from 2022-Q3).
function Set-TextValue($ws, $addr, $val) {
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $val
}

Set-TextValue $q4Ws "D2" "0.23"
Set-TextValue $q4Ws "E2" "94.47"
Set-TextValue $q4Ws "F2" "3.38"
Set-TextValue $q4Ws "G2" "0.0078"

Set-TextValue $q4Ws "D3" "0.16"
Set-TextValue $q4Ws "E3" "93.15"
Set-TextValue $q4Ws "F3" "1.23"
Set-TextValue $q4Ws "G3" "0.0020"
$q4Ws.Range("H3").Value = 9

Set-TextValue $q4Ws "E4" "93.15"
Set-TextValue $q4Ws "F4" "1.23"
$q4Ws.Range("H4").Value = 9

# --- 2. Update the "总计" summary sheet: shift every existing quarter's
#        row down by one and insert the new "2022-Q4" row at the top. ---
$totalWs.Range("B2").Value = "2022-Q4"
$totalWs.Range("C2").Value = 3
$totalWs.Range("D2").Value = 0.01

$totalWs.Range("B3").Value = "2022-Q3"
$totalWs.Range("C3").Value = 3
$totalWs.Range("D3").Value = 0.01

$totalWs.Range("B4").Value = "2022-Q2"
$totalWs.Range("C4").Value = 4
$totalWs.Range("D4").Value = 0.15

$totalWs.Range("B5").Value = "2022-Q1"
$totalWs.Range("C5").Value = 4
$totalWs.Range("D5").Value = 0.22

$totalWs.Range("B6").Value = "2021-Q4"
$totalWs.Range("C6").Value = 3
$totalWs.Range("D6").Value = 0.01

$totalWs.Range("B7").Value = "2021-Q3"
$totalWs.Range("C7").Value = 2
$totalWs.Range("D7").Value = 0.02

$totalWs.Range("A8").Value = 6
$totalWs.Range("B8").Value = "2021-Q2"
$totalWs.Range("C8").Value = 5
$totalWs.Range("D8").Value = 0.11

# Copy A2's formatting (bold / border / centered) onto the freshly added
# A8 cell so the new summary row matches the rest of the table.
$totalWs.Range("A2").Copy()
$totalWs.Range("A8").PasteSpecial(-4122)

# --- 3. Restore the originally active tab (the oldest quarter, 2021-Q2,
#        was the selected sheet before the edit). ------------------------
$lastWs = $wb.Worksheets.Item($wb.Worksheets.Count)
$lastWs.Activate()
$totalWs.Range("A1").Select()
